$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "2106"
$ws.Range("E17").Value = "2105"
$ws.Range("E18").Value = "2104"
$ws.Range("E19").Value = "2103"
$ws.Range("E20").Value = "2102"
$ws.Range("E21").Value = "2101"
$ws.Range("E22").Value = "2012"
$ws.Range("E23").Value = "2011"
$ws.Range("E24").Value = "2010"

$ws.Range("F16").Value = 34874
$ws.Range("F24").Value = 49820
